$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the "Delaying Practicum III..." sentence was split across three
# runs (with w:proofErr gramStart/gramEnd bracketing "also has"). Collapse it
# into a single run with no grammar-check markers, leaving the preceding
# space-run untouched.
# ---------------------------------------------------------------------------
$target = "Delaying Practicum III completion by a semester also has thrown my initial planning for a loop."
$full = $d.Content.Text
$idx = $full.IndexOf($target)
if ($idx -ge 0) {
    $oldRange = $d.Range($idx, $idx + $target.Length)
    [void]$oldRange.Delete()
    $insertionPoint = $d.Range($idx, $idx)
    [void]$insertionPoint.InsertAfter($target)
}

# ---------------------------------------------------------------------------
# Change 2: add a new 2/24/2024 journal entry (Niagara system rework) right
# after the 11/30/2023 entry, preserving the trailing blank paragraph that
# already exists before the section break.
# ---------------------------------------------------------------------------
$apos = [char]0x2019
$niagaraText = "I replaced the constant Niagara system effect within the "
$niagaraTail = " with the ability for the character to spawn in the effect when the press the key associated with power activation. I" + $apos + "m hoping this leads to a less-intensive and less confusing result. It" + $apos + "s location does need to be constantly updated, which I have not done yet."

$newParasXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:r><w:t>2/24/2024</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:r><w:tab/><w:t xml:space='preserve'>${niagaraText}</w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>FP_Character</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t xml:space='preserve'>${niagaraTail}</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
"@

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$anchor = $lastPara.Range.Start
$insertPoint = $d.Range($anchor, $anchor)
[void]$insertPoint.InsertXML($newParasXml)
